$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 244
$ws.Range("I5").Value = 244
$ws.Range("K5").Value = 244
$ws.Range("M5").Value = -129

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1822.1111
$ws.Range("I17").Value = 1549
$ws.Range("J17").Value = 1856.25
$ws.Range("K17").Value = 4647
$ws.Range("L17").Value = 5568.75
$ws.Range("M17").Value = -4479
$ws.Range("N17").Value = -5904.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 3000
$ws.Range("N86").Value = -5246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 15000
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3995.5
$ws.Range("I106").Value = 3995.5
$ws.Range("K106").Value = 3995.5
$ws.Range("M106").Value = -3364.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2990.6667
$ws.Range("I2").Value = 1761.25
$ws.Range("K2").Value = 1761.25
$ws.Range("M2").Value = -1648.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6860.75
$ws.Range("I32").Value = 6120.909
$ws.Range("K32").Value = 6120.909
$ws.Range("M32").Value = -5833.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 28800
$ws.Range("J76").Value = 28800
$ws.Range("L76").Value = 28800
$ws.Range("N76").Value = -29476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 28800
$ws.Range("J79").Value = 28800
$ws.Range("L79").Value = 28800
$ws.Range("N79").Value = -31140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 913
$ws.Range("I110").Value = 619.75
$ws.Range("K110").Value = 619.75
$ws.Range("M110").Value = 1425.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2990.6667
$ws.Range("I116").Value = 1761.25
$ws.Range("K116").Value = 1761.25
$ws.Range("M116").Value = 532.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2990.6667
$ws.Range("I3").Value = 1761.25
$ws.Range("K3").Value = 1761.25
$ws.Range("M3").Value = -1647.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 1000
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1000
$ws.Range("L24").ClearContents()
$ws.Range("M24").Value = -765
$ws.Range("N24").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6666.3335
$ws.Range("I105").Value = 5999.5
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 5999.5
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -4252.5
$ws.Range("N105").Value = -11494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4198.8
$ws.Range("I107").Value = 2748.5
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 2748.5
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -828.5
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 54
$ws.Range("I2").Value = 54
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2666.625
$ws.Range("I31").Value = 889
$ws.Range("K31").Value = 889
$ws.Range("M31").Value = -594

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2666.625
$ws.Range("I34").Value = 889
$ws.Range("K34").Value = 889
$ws.Range("M34").Value = -687

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1950
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1200383.8
$ws.Range("I4").Value = 1250480
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 3751440
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -3751328
$ws.Range("N4").Value = -3000221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 461.5
$ws.Range("J107").Value = 461.5
$ws.Range("L107").Value = 1384.5
$ws.Range("N107").Value = -5224.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3348
$ws.Range("J131").Value = 3644.2222
$ws.Range("L131").Value = 10932.6666
$ws.Range("N131").Value = -21012.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1997.5
$ws.Range("J132").Value = 1995
$ws.Range("L132").Value = 17955
$ws.Range("N132").Value = -23015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4233.1665
$ws.Range("I80").Value = 3833
$ws.Range("J80").Value = 4633.3335
$ws.Range("K80").Value = 3833
$ws.Range("L80").Value = 4633.3335
$ws.Range("M80").Value = -2835
$ws.Range("N80").Value = -6629.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4233.1665
$ws.Range("I83").Value = 3833
$ws.Range("J83").Value = 4633.3335
$ws.Range("K83").Value = 19165
$ws.Range("L83").Value = 23166.6675
$ws.Range("M83").Value = -14173
$ws.Range("N83").Value = -33150.6675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2149
$ws.Range("I102").Value = 2149
$ws.Range("K102").Value = 2149
$ws.Range("M102").Value = -527

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4883.1665
$ws.Range("I113").Value = 4883.1665
$ws.Range("K113").Value = 4883.1665
$ws.Range("M113").Value = -2713.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1049.8
$ws.Range("I122").Value = 966.3333
$ws.Range("J122").Value = 1175
$ws.Range("K122").Value = 2898.9999
$ws.Range("L122").Value = 3525
$ws.Range("M122").Value = -448.9998999999998
$ws.Range("N122").Value = -8425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2400
$ws.Range("I16").Value = 2400
$ws.Range("K16").Value = 2400
$ws.Range("M16").Value = -2230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 7511
$ws.Range("I19").Value = 7511
$ws.Range("K19").Value = 7511
$ws.Range("M19").Value = -7341

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 4000
$ws.Range("N22").Value = -4590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 4000
$ws.Range("N27").Value = -4214

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3341.111
$ws.Range("I61").Value = 2414.8
$ws.Range("K61").Value = 2414.8
$ws.Range("M61").Value = -2212.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4216.5
$ws.Range("I100").Value = 3059.8
$ws.Range("K100").Value = 3059.8
$ws.Range("M100").Value = -2518.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3341.111
$ws.Range("I113").Value = 2414.8
$ws.Range("K113").Value = 2414.8
$ws.Range("M113").Value = -244.8000000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4999.25
$ws.Range("I122").Value = 4999.25
$ws.Range("K122").Value = 14997.75
$ws.Range("M122").Value = -12547.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3489.5
$ws.Range("I107").Value = 1734.25
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 5202.75
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = -3282.75
$ws.Range("N107").Value = -24840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5688
$ws.Range("I122").Value = 5416.5
$ws.Range("J122").Value = 6502.5
$ws.Range("K122").Value = 16249.5
$ws.Range("L122").Value = 19507.5
$ws.Range("M122").Value = -13799.5
$ws.Range("N122").Value = -24407.5
